$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reword the "current top-5 goal scorer" question in A16
$ws.Range("A16").Value = "Name a current player of the team that is ranked among the Top 5 striker of all time?"

# 2. Highlight the question/answer block (rows 10-18) with a yellow fill,
#    matching the range that was selected and formatted in Excel.
$yellow = 65535  # RGB(255,255,0) -> BGR COM color value
$ws.Range("A10:B18").Interior.Color = $yellow
$ws.Range("C14").Interior.Color = $yellow

# 3. Scroll the view down and move the active selection to A19,
#    as it was left after the edit.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A19").Select() | Out-Null
